# HYP_VMH_Selected.xlsx - tidy/summarize header sample IDs
# (commit: "tidying and summarizzation of VMH / ARC smfish cell phenotypes
#  (gene+/-, gene +/-, ...) by sex")
#
# The underlying data (counts per phenotype row) is unchanged; only three
# of the sample-ID column headers in row 1 were corrected/relabeled:
#   C1: Br1223_VMH_ARC -> Br1225_VMH_ARC
#   F1: Br8471_VMH_ARC -> Br8741_VMH_ARC
#   I1: VMH_ARC         -> Br8667_VMH_ARC

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Br1225_VMH_ARC"
$ws.Range("F1").Value = "Br8741_VMH_ARC"
$ws.Range("I1").Value = "Br8667_VMH_ARC"
